$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work in manual calculation mode so that touching one cell does not
# trigger a full, clock-dependent recalculation of every TODAY()-based
# "calculate always" formula in the sheet (B/H/J/K columns). Each
# formula cell we need refreshed is re-entered explicitly below, which
# (like real Excel) computes just that cell on entry.
$excel.Calculation = -4135   # xlCalculationManual

# ---- Real data edits (column C, the daily new-case counts) ----
$ws.Range("C251").Value2 = 851
$ws.Range("C281").Value2 = 120
$ws.Range("C290").Value2 = 122
$ws.Range("C291").Value2 = 72
$ws.Range("C293").Value2 = 116
$ws.Range("C294").Value2 = 75
$ws.Range("C295").Value2 = 16

# ---- New data for row 295 (E, F, G, L, M) ----
$ws.Range("E295").Value2 = 19
$ws.Range("F295").Value2 = 15
$ws.Range("G295").Value2 = 85

# L295/M295 are formatted as Text (@); round-trip the number format so the
# literal 0 is stored as a genuine number (matching every other cell in
# those columns) instead of being coerced into a text value.
$ws.Range("L295").NumberFormat = "General"
$ws.Range("L295").Value2 = 0
$ws.Range("L295").NumberFormat = "@"

$ws.Range("M295").NumberFormat = "General"
$ws.Range("M295").Value2 = 0
$ws.Range("M295").NumberFormat = "@"

# ---- Refresh the cumulative-sum formula chain in column B (si=19 group) ----
# B251:B259 belong to the B196:B259 shared group; B260:B295 (and the
# newly-populated B296) belong to the B260:B310 shared group. Re-entering
# each cell's formula recalculates just that cell, in row order, so the
# running total ripples correctly without touching B297 onward.
for ($r = 251; $r -le 296; $r++) {
    $prev = $r - 1
    $ws.Range("B$r").Formula = "=IF(TODAY()>A$prev,B$prev+C$r,`"`")"
}

# ---- Refresh column H (G+E) for the rows whose inputs changed ----
foreach ($r in 295, 296) {
    $prev = $r - 1
    $ws.Range("H$r").Formula = "=IF(TODAY()>A$prev,G$r+E$r,`"`")"
}

# ---- Refresh column J (running total) and K (L+M) for row 296, now in range ----
foreach ($r in 296) {
    $prev = $r - 1
    $ws.Range("J$r").Formula = "=IF(TODAY()>A$prev,J$prev+K$r,`"`")"
    $ws.Range("K$r").Formula = "=IF(TODAY()>A$prev,L$r+M$r,`"`")"
}
